$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.07
$ws.Range("K2").Value = 9
$ws.Range("N2").Value = 2.1
$ws.Range("O2").Value = 1.7

# Row 6
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 3.15
$ws.Range("I6").Value = 3.7
$ws.Range("L6").Value = 1.47
$ws.Range("M6").Value = 2.35
$ws.Range("N6").Value = 2.32
$ws.Range("O6").Value = 1.47
$ws.Range("P6").Value = 1.5
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 2.07
$ws.Range("S6").Value = 1.6
$ws.Range("T6").Value = 5.5
$ws.Range("U6").Value = 8
$ws.Range("V6").Value = 9.25
$ws.Range("W6").Value = 17
$ws.Range("X6").Value = 19.5
$ws.Range("Y6").Value = 40
$ws.Range("Z6").Value = 6.9
$ws.Range("AA6").Value = 6.3
$ws.Range("AB6").Value = 19.5
$ws.Range("AC6").Value = 120
$ws.Range("AE6").Value = 8.25
$ws.Range("AF6").Value = 18
$ws.Range("AG6").Value = 13.5
$ws.Range("AH6").Value = 55
$ws.Range("AI6").Value = 45
$ws.Range("AJ6").Value = 60

# Row 8
$ws.Range("G8").Value = 1.72
$ws.Range("H8").Value = 3.9
$ws.Range("I8").Value = 4
$ws.Range("R8").Value = 1.65
$ws.Range("S8").Value = 2
$ws.Range("X8").Value = 13
$ws.Range("Z8").Value = 13.5
$ws.Range("AH8").Value = 60

# Row 9
$ws.Range("H9").Value = 3.65
$ws.Range("I9").Value = 3.2
$ws.Range("L9").Value = 1.19
$ws.Range("M9").Value = 3.75
$ws.Range("N9").Value = 1.57
$ws.Range("O9").Value = 2.1
$ws.Range("R9").Value = 1.52
$ws.Range("S9").Value = 2.2
$ws.Range("T9").Value = 9.75
$ws.Range("X9").Value = 14
$ws.Range("Y9").Value = 20
$ws.Range("Z9").Value = 13.5
$ws.Range("AA9").Value = 7.4
$ws.Range("AD9").Value = 250
$ws.Range("AF9").Value = 19.5
$ws.Range("AI9").Value = 25
$ws.Range("AJ9").Value = 28

# Row 14
$ws.Range("G14").Value = 3.35
$ws.Range("H14").Value = 3.2
$ws.Range("I14").Value = 2.02
$ws.Range("N14").Value = 2.15
$ws.Range("O14").Value = 1.55
$ws.Range("T14").Value = 7.1
$ws.Range("U14").Value = 13
$ws.Range("V14").Value = 10.25
$ws.Range("W14").Value = 35
$ws.Range("X14").Value = 26
$ws.Range("Y14").Value = 37
$ws.Range("Z14").Value = 7.7
$ws.Range("AA14").Value = 5.5
$ws.Range("AB14").Value = 14
$ws.Range("AC14").Value = 70
$ws.Range("AE14").Value = 5.3
$ws.Range("AF14").Value = 7.3
$ws.Range("AG14").Value = 7.6
$ws.Range("AH14").Value = 14
$ws.Range("AI14").Value = 15

# Row 17
$ws.Range("T17").Value = 7.2
$ws.Range("U17").Value = 10

# Row 18
$ws.Range("I18").Value = 3.8

# Row 35
$ws.Range("N35").Value = 1.62
$ws.Range("O35").Value = 2.25
